$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on price cells whose new values look like plain numbers,
# so Excel keeps them as text (matching the original inline-string formatting)
# instead of silently converting to a numeric value.
$textCells = @("D5","D6","D9","D10","D11","D13","D15","D17","D21","D22","D23","D24","D27","D29","D30","D32","D35","D36","D37","D40","D42","D45","D46","D48")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated price / 1h-volume figures scraped by the Action run.
$ws.Range("D2").Value = "42.941.42"
$ws.Range("E2").Value = "  -1.06%  "
$ws.Range("D3").Value = "2.308.84"
$ws.Range("E3").Value = "  -0.03%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").Value = "302.42"
$ws.Range("E5").Value = "  -1.98%  "
$ws.Range("D6").Value = "100.01"
$ws.Range("E6").Value = "  -4.92%  "
$ws.Range("E7").Value = "  -3.85%  "
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("D9").Value = "0.503"
$ws.Range("E9").Value = "  -2.58%  "
$ws.Range("D10").Value = "34.80"
$ws.Range("E10").Value = "  -3.04%  "
$ws.Range("D11").Value = "0.0792"
$ws.Range("E11").Value = "  -2.19%  "
$ws.Range("E12").Value = "  +0.45%  "
$ws.Range("D13").Value = "6.71"
$ws.Range("E13").Value = "  -3.41%  "
$ws.Range("D14").Value = "2.665.25"
$ws.Range("E14").Value = "  -0.16%  "
$ws.Range("D15").Value = "15.54"
$ws.Range("E15").Value = "  +2.97%  "
$ws.Range("D16").Value = "2.273.85"
$ws.Range("E16").Value = "  -1.60%  "
$ws.Range("D17").Value = "0.795"
$ws.Range("E17").Value = "  -0.65%  "
$ws.Range("D18").Value = "42.823.48"
$ws.Range("E18").Value = "  -1.21%  "
$ws.Range("E19").Value = "  -1.23%  "
$ws.Range("E20").Value = "  -1.67%  "
$ws.Range("D21").Value = "6.04"
$ws.Range("E21").Value = "  -2.48%  "
$ws.Range("D22").Value = "67.83"
$ws.Range("E22").Value = "  -0.16%  "
$ws.Range("D23").Value = "235.79"
$ws.Range("E23").Value = "  -1.85%  "
$ws.Range("D24").Value = "1.95"
$ws.Range("E24").Value = "  -3.68%  "
$ws.Range("E25").Value = "  -3.66%  "
$ws.Range("E26").Value = "  -0.08%  "
$ws.Range("D27").Value = "24.74"
$ws.Range("E27").Value = "  -0.97%  "
$ws.Range("E28").Value = "  -6.27%  "
$ws.Range("D29").Value = "34.51"
$ws.Range("E29").Value = "  -4.78%  "
$ws.Range("D30").Value = "164.40"
$ws.Range("E30").Value = "  +0.87%  "
$ws.Range("E31").Value = "  -4.93%  "
$ws.Range("D32").Value = "0.999"
$ws.Range("E32").Value = "  -0.09%  "
$ws.Range("E33").Value = "  -4.18%  "
$ws.Range("E34").Value = "  -4.50%  "
$ws.Range("D35").Value = "4.51"
$ws.Range("E35").Value = "  -1.76%  "
$ws.Range("D36").Value = "16.72"
$ws.Range("E36").Value = "  -8.18%  "
$ws.Range("D37").Value = "0.0696"
$ws.Range("E37").Value = "  -4.78%  "
$ws.Range("E38").Value = "  -3.94%  "
$ws.Range("E39").Value = "  -3.20%  "
$ws.Range("D40").Value = "0.100"
$ws.Range("E40").Value = "  -4.71%  "
$ws.Range("E41").Value = "  -3.61%  "
$ws.Range("D42").Value = "2.50"
$ws.Range("E42").Value = "  +0.13%  "
$ws.Range("D43").Value = "1.974.26"
$ws.Range("E43").Value = "  +0.54%  "
$ws.Range("E44").Value = "  -3.07%  "
$ws.Range("D45").Value = "18.36"
$ws.Range("E45").Value = "  -2.16%  "
$ws.Range("D46").Value = "10.23"
$ws.Range("E46").Value = "  -0.34%  "
$ws.Range("E47").Value = "  -5.34%  "
$ws.Range("D48").Value = "55.49"
$ws.Range("E48").Value = "  -3.85%  "
$ws.Range("D49").Value = "2.531.59"
$ws.Range("E49").Value = "  -0.16%  "
$ws.Range("E50").Value = "  -4.40%  "
$ws.Range("E51").Value = "  +0.27%  "
